$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.621.68'
$ws.Cells.Item(2, 5).Value = '  -1.49%  '
$ws.Cells.Item(3, 4).Value = '2.630.74'
$ws.Cells.Item(3, 5).Value = '  +0.89%  '
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(5, 4).Value = '536.35'
$ws.Cells.Item(5, 5).Value = '  -0.11%  '
$ws.Cells.Item(6, 4).Value = '142.92'
$ws.Cells.Item(6, 5).Value = '  +0.98%  '
$ws.Cells.Item(7, 5).Value = '  -0.09%  '
$ws.Cells.Item(8, 5).Value = '  -0.18%  '
$ws.Cells.Item(9, 4).Value = '2.638.59'
$ws.Cells.Item(9, 5).Value = '  +0.87%  '
$ws.Cells.Item(10, 5).Value = '  +8.17%  '
$ws.Cells.Item(11, 5).Value = '  -1.71%  '
$ws.Cells.Item(12, 5).Value = '  -0.21%  '
$ws.Cells.Item(13, 5).Value = '  +0.93%  '
$ws.Cells.Item(14, 4).Value = '3.098.67'
$ws.Cells.Item(15, 4).Value = '58.565.90'
$ws.Cells.Item(15, 5).Value = '  -1.46%  '
$ws.Cells.Item(16, 5).Value = '  +1.07%  '
$ws.Cells.Item(17, 4).Value = '2.636.93'
$ws.Cells.Item(17, 5).Value = '  +0.55%  '
$ws.Cells.Item(18, 5).Value = '  -0.85%  '
$ws.Cells.Item(19, 5).Value = '  +1.16%  '
$ws.Cells.Item(20, 4).Value = '334.66'
$ws.Cells.Item(20, 5).Value = '  -1.98%  '
$ws.Cells.Item(21, 5).Value = '  +0.49%  '
$ws.Cells.Item(22, 4).Value = '6.25'
$ws.Cells.Item(22, 5).Value = '  -2.02%  '
$ws.Cells.Item(23, 5).Value = '  +0.01%  '
$ws.Cells.Item(24, 4).Value = '66.34'
$ws.Cells.Item(24, 5).Value = '  -1.71%  '
$ws.Cells.Item(25, 5).Value = '  +1.37%  '
$ws.Cells.Item(26, 5).Value = '  -0.96%  '
$ws.Cells.Item(27, 5).Value = '  +0.07%  '
$ws.Cells.Item(28, 4).Value = '7.15'
$ws.Cells.Item(28, 5).Value = '  -1.25%  '
$ws.Cells.Item(29, 5).Value = '  -1.16%  '
$ws.Cells.Item(30, 5).Value = '  -0.04%  '
$ws.Cells.Item(31, 5).Value = '  -1.19%  '
$ws.Cells.Item(32, 4).Value = '5.83'
$ws.Cells.Item(32, 5).Value = '  +0.07%  '
$ws.Cells.Item(33, 4).Value = '18.75'
$ws.Cells.Item(33, 5).Value = '  -0.63%  '
$ws.Cells.Item(34, 4).Value = '150.42'
$ws.Cells.Item(34, 5).Value = '  +0.57%  '
$ws.Cells.Item(35, 5).Value = '  -1.98%  '
$ws.Cells.Item(36, 4).Value = '37.15'
$ws.Cells.Item(36, 5).Value = '  -0.10%  '
$ws.Cells.Item(37, 4).Value = '1.10'
$ws.Cells.Item(37, 5).Value = '  -0.97%  '
$ws.Cells.Item(38, 4).Value = '0.844'
$ws.Cells.Item(38, 5).Value = '  +0.97%  '
$ws.Cells.Item(39, 5).Value = '  -2.98%  '
$ws.Cells.Item(40, 5).Value = '  -1.77%  '
$ws.Cells.Item(41, 5).Value = '  +1.05%  '
$ws.Cells.Item(42, 4).Value = '282.13'
$ws.Cells.Item(42, 5).Value = '  +2.85%  '
$ws.Cells.Item(43, 5).Value = '  -0.10%  '
$ws.Cells.Item(44, 4).Value = '0.602'
$ws.Cells.Item(44, 5).Value = '  +0.35%  '
$ws.Cells.Item(45, 4).Value = '10.68'
$ws.Cells.Item(45, 5).Value = '  -0.46%  '
$ws.Cells.Item(46, 4).Value = '19.12'
$ws.Cells.Item(46, 5).Value = '  +3.05%  '
$ws.Cells.Item(47, 4).Value = '0.0532'
$ws.Cells.Item(47, 5).Value = '  +1.61%  '
$ws.Cells.Item(48, 4).Value = '0.0937'
$ws.Cells.Item(48, 5).Value = '  -1.88%  '
$ws.Cells.Item(49, 4).Value = '0.0225'
$ws.Cells.Item(49, 5).Value = '  +0.72%  '
$ws.Cells.Item(50, 4).Value = '1.947.42'
$ws.Cells.Item(50, 5).Value = '  -0.04%  '
$ws.Cells.Item(51, 4).Value = '4.45'
$ws.Cells.Item(51, 5).Value = '  -1.27%'
